# Add "recurring" (boolean) and "recurrence" (text) fields to the bulk-upload
# template, mirroring the existing expense/income rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G1").Value = "recurring"
$ws.Range("H1").Value = "recurrence"

# Row 2 - recurring weekly
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = "weekly"

# Row 3
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = "monthly"

# Row 4 - not recurring
$ws.Range("G4").Value = $false

# Row 5 - not recurring
$ws.Range("G5").Value = $false

# Row 6 - not recurring
$ws.Range("G6").Value = $false

# Row 7 - recurring daily
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = "daily"

# Row 8 - not recurring
$ws.Range("G8").Value = $false

# Row 9 - not recurring
$ws.Range("G9").Value = $false

# Match the column sizing used for the new "recurring" column.
$ws.Columns.Item(7).ColumnWidth = 17

# Restore the selection to the cell last edited, as captured in the source file.
$ws.Range("H7").Select()
